$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-04-30"

# Update the header label for April in column A, row 5
$ws.Range("A5").Value = "April (through 04-30)"

# Update April (row 5) counts
$ws.Range("B5").Value = 23
$ws.Range("C5").Value = 34
$ws.Range("D5").Value = 64
$ws.Range("E5").Value = 49
$ws.Range("F5").Value = 45
$ws.Range("G5").Value = 64
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 117

# Update Total (row 6) counts
$ws.Range("B6").Value = 89
$ws.Range("C6").Value = 162
$ws.Range("D6").Value = 253
$ws.Range("E6").Value = 246
$ws.Range("F6").Value = 155
$ws.Range("G6").Value = 262
$ws.Range("H6").Value = 523
$ws.Range("I6").Value = 552
